# Apply commit: "[ADDITIONAL SCRAPING] added scraping code for extra browling
# attributes and excel sheets"
#
# 1) Trim stray empty placeholder cells out of "ODI Batting Extra" (they used
#    to be written as empty <c t="inlineStr"/> cells; the re-scrape only keeps
#    cells that actually carry data).
# 2) Add a brand-new "ODI Bowling Extra" sheet (sheetId 5) with the bowling
#    extras scraped from each match: MAIDEN_OVERS and PERCENT_WICKETS_OF_ALL.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "ODI Batting Extra" - clear out the now-empty placeholder cells
# ---------------------------------------------------------------------------
$batExtra = $wb.Worksheets.Item("ODI Batting Extra")

# Rows whose B:E block has no data at all (only the MAN_OF_MATCH / F column
# still holds a value).
$clearBtoE = @(3, 4, 7, 8, 10, 13, 14)
foreach ($r in $clearBtoE) {
    $batExtra.Range("B$r`:E$r").ClearContents()
}

# Row 2 keeps its BATTING_POSITION (B2) and MAN_OF_MATCH (F2) values; only
# C2:E2 are empty.
$batExtra.Range("C2:E2").ClearContents()

# Rows 5 and 11 only have the PERCENT_RUNS_OF_TOTAL (E) column empty.
$batExtra.Range("E5").ClearContents()
$batExtra.Range("E11").ClearContents()

# Rows 15-21 have no data at all beyond the MATCH_CODE (A) column.
$clearBtoF = @(15, 16, 17, 18, 19, 20, 21)
foreach ($r in $clearBtoF) {
    $batExtra.Range("B$r`:F$r").ClearContents()
}

# ---------------------------------------------------------------------------
# 2) Add the "ODI Bowling Extra" sheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlExtra.Name = "ODI Bowling Extra"

# Header row - reuse the bold/bordered/centered header style already used by
# the other "* Extra" sheets.
$bowlExtra.Cells.Item(1, 1).Value = "MATCH_CODE"
$bowlExtra.Cells.Item(1, 2).Value = "MAIDEN_OVERS"
$bowlExtra.Cells.Item(1, 3).Value = "PERCENT_WICKETS_OF_ALL"

$batExtra.Range("A1").Copy()
$bowlExtra.Range("A1").PasteSpecial(-4122)
$batExtra.Range("B1").Copy()
$bowlExtra.Range("B1").PasteSpecial(-4122)
$batExtra.Range("C1").Copy()
$bowlExtra.Range("C1").PasteSpecial(-4122)

# Data rows, one per match, matching "ODI Batting Extra"'s MATCH_CODE order.
# Every value is written as text (leading "'" forces text storage instead of
# Excel auto-coercing numeric-looking strings / percentages into numbers),
# matching the original scrape where every cell is a string.
$rows = @(
    @("4210", "0", "10.00%"),
    @("4211", "0", ""),
    @("4231", "0", "20.00%"),
    @("4232", "0", "10.00%"),
    @("4233", "", ""),
    @("4302", "", ""),
    @("4309", "1", "40.00%"),
    @("4322", "0", ""),
    @("4331", "1", "10.00%"),
    @("4356", "", ""),
    @("4357", "", ""),
    @("4375", "0", ""),
    @("4376", "1", "20.00%"),
    @("4413", "0", "10.00%"),
    @("4414", "", ""),
    @("4417", "1", ""),
    @("4449", "1", ""),
    @("4450", "", ""),
    @("4521", "1", "20.00%"),
    @("4523", "", "")
)

$r = 2
foreach ($row in $rows) {
    $bowlExtra.Cells.Item($r, 1).Value = "'" + $row[0]
    $bowlExtra.Cells.Item($r, 2).Value = "'" + $row[1]
    $bowlExtra.Cells.Item($r, 3).Value = "'" + $row[2]
    $r = $r + 1
}
